$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# out_vars: append row 12 with the June 11 (serial 43993) summary,
# copying row 11's look (date format in col A, general+wrap B:J).
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("out_vars")

$ws1.Range("A11").Copy()
$ws1.Range("A12").PasteSpecial(-4122)
$ws1.Range("B11:J11").Copy()
$ws1.Range("B12:J12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A12").Value = 43993
$ws1.Range("B12").Value = 133974
$ws1.Range("C12").Value = 191465
$ws1.Range("D12").Value = 55700
$ws1.Range("E12").Value = 15944
$ws1.Range("F12").Value = 33.006404227685969
$ws1.Range("G12").Value = 44220
$ws1.Range("H12").Value = 4087
$ws1.Range("I12").Value = 4256
$ws1.Range("J12").Value = 381139

$ws1.Range("A12").Select()

# -----------------------------------------------------------------
# dates_dx: row 12 placeholders already carry their style - just
# fill in the June 11 counts.
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dates_dx")
$ws2.Range("A12").Value = 43993
$ws2.Range("B12").Value = 0
$ws2.Range("C12").Value = 1
$ws2.Range("D12").Value = 1
$ws2.Range("E12").Value = 1
$ws2.Range("F12").Value = 0
$ws2.Range("G12").Value = 0
$ws2.Range("H12").Value = 0
$ws2.Range("I12").Value = 4

$ws2.Range("J12").Select()

# -----------------------------------------------------------------
# dates_sx: append row 12 with data and grow the trailing blank
# placeholder row to row 13 (matching col A's style elsewhere).
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("dates_sx")

$ws3.Range("A11").Copy()
$ws3.Range("A12").PasteSpecial(-4122)
$ws3.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("A12").Value = 43993
$ws3.Range("B12").Value = 0
$ws3.Range("C12").Value = 1
$ws3.Range("D12").Value = 0
$ws3.Range("E12").Value = 1
$ws3.Range("F12").Value = 1
$ws3.Range("G12").Value = 1
$ws3.Range("H12").Value = 0
$ws3.Range("I12").Value = 1
$ws3.Range("J12").Value = 1
$ws3.Range("K12").Value = 0
$ws3.Range("L12").Value = 0

$ws3.Range("C20").Select()

# -----------------------------------------------------------------
# dates_deaths: row 12 already exists (blank) - fill in the values
# and pick up col A's date style (was a generic blank-row style).
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("dates_deaths")

$ws4.Range("A11").Copy()
$ws4.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Range("A12").Value = 43993
$ws4.Range("B12").Value = 1
$ws4.Range("C12").Value = 0
$ws4.Range("D12").Value = 2
$ws4.Range("E12").Value = 1
$ws4.Range("F12").Value = 1
$ws4.Range("G12").Value = 2
$ws4.Range("H12").Value = 2

$ws4.Range("I12").Select()

# -----------------------------------------------------------------
# control_obs: fill in the June 11 (column L) figures and extend
# the running-total formula from column K into column L.
# -----------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("control_obs")

$ws5.Range("L1").Value = 43993
$ws5.Range("L2").Value = 3369
$ws5.Range("L3").Value = 3184
$ws5.Range("L4").Value = 3184
$ws5.Range("L5").Value = 3184
$ws5.Range("L6").Value = 3184
$ws5.Range("L7").Value = 2481
$ws5.Range("L8").Value = 5006
$ws5.Range("L10").Value = 151
$ws5.Range("L11").Value = 151
$ws5.Range("L12").Value = 151
$ws5.Range("L13").Value = 151
$ws5.Range("L14").Value = 151
$ws5.Range("L15").Value = 128
$ws5.Range("L16").Value = 163
$ws5.Range("L18").Value = 813
$ws5.Range("L20").Formula = "=SUM(L2:L18)"

$ws5.Range("L25").Select()
$ws5.Activate()

Write-Output "edit applied"
